$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.623.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.443.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.72%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'575.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.85%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.62%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.72%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.439.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.99%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.61%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.59%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.26%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.372.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.440.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.87%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'328.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.41%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.35%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.21%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'634.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +7.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0968"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.62%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.23%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.88%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.27%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.19%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.41%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'18.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.52%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'146.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.56%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'42.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.89%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'145.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.05%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.50%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'19.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.21%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.30%  "
$ws.Range("E51").Style = "Normal"
Write-Output "Applied crypto price/volume updates"
